# Update the PortageII API comparison sheet (Sheet1) to cover the 4.0 release:
#  - B1 / B11: section header titles bumped from "...to PortageII 3.0[.../(p.2)]"
#              to "...to PortageII 4.0[.../(p.2)]"
#  - B4: getVersion() API row now documents both the 3.0 and 4.0 behaviour
#  - B21 / B22: incremental-training API rows now note the 3.0 "n/a" baseline
#              alongside the existing 4.0 method names
# Finally, move the saved cell selection to B3 (where the reviewers were
# pointed while looking the sheet over).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$NL = [char]10

$ws.Range("B1").Value = "Evolution of the PortageII API from PortageII 1.0 to PortageII 4.0"

$ws.Range("B4").Value = "3.0: getVersion() exists and returns ""PortageII-3.0""; WSDL says PortageII 3.0 (2016)" + $NL + "4.0: getVersion() exists and returns ""PortageII-4.0""; WSDL says PortageII 4.0 (2018)"

$ws.Range("B11").Value = "Evolution of the PortageII API from PortageII 1.0 to PortageII 4.0 (p. 2)"

$ws.Range("B21").Value = "3.0: n/a" + $NL + "4.0: incrAddSentence()"

$ws.Range("B22").Value = "3.0: n/a" + $NL + "4.0: incrStatus()"

$ws.Activate()
$ws.Range("B3").Select()
